$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) column C holds a date serial number for rows 2-74.
# Update every occurrence of 45186 to 45188 (2023-09-17 -> 2023-09-19).
for ($r = 2; $r -le 74; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
